$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

$ws.Range("D8").Value = "M.2 E connector"
$ws.Range("E8").Value = "SparkFun"

$ws.Hyperlinks.Add($ws.Range("D8"), "https://www.sparkfun.com/products/15427", "", "", "M.2 E connector")
$ws.Range("D8").Style = "Hyperlink"

$wb.Save()
